$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----- Column widths -----
# Raw OOXML width = ColumnWidth + ~0.7142857 (MDW padding) for this Tahoma-based workbook.
# Column A -> raw width 7 ; Column B -> raw width ~8.375 (closest achievable via pixel rounding)
$ws.Columns.Item(1).ColumnWidth = 6.285714285714286
$ws.Columns.Item(2).ColumnWidth = 7.645

# ----- Enter all text values first (keeps shared-string table order close to the
#       original authoring order: title, ID, Page, Detail, 01..21) -----
$ws.Range("A1").Value = "Error Tabel"
$ws.Range("A2").Value = "ID "
$ws.Range("B2").Value = "Page"
$ws.Range("C2").Value = "Detail"

$ids = @("01","02","03","04","05","06","07","08","09","10","11","12","13","14","15","16","17","18","19","20","21")
for ($i = 0; $i -lt $ids.Count; $i++) {
    $row = 3 + $i
    $ws.Range("A$row").Value = $ids[$i]
}

# ----- Formatting -----
# Row 2 headers "Page"/"Detail" centered
$ws.Range("B2:C2").HorizontalAlignment = -4108

# Title row centered and merged across A1:C1
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1:C1").Merge() | Out-Null

# Column A values (the "ID " header and "01".."21") stored as centered text so the
# leading zero is preserved
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").NumberFormat = "@"
for ($i = 0; $i -lt $ids.Count; $i++) {
    $row = 3 + $i
    $cell = $ws.Range("A$row")
    $cell.HorizontalAlignment = -4108
    $cell.NumberFormat = "@"
    $cell.Value = $ids[$i]
}

# ----- Selection shown in the saved view -----
$ws.Range("B3").Select() | Out-Null
